$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label fixes (swap to correct variable names in line_of_thrust_center)
$ws.Range("V1").Value = "a_dx"
$ws.Range("W1").Value = "a_dy"

# Recalculated numeric results (ripple from the variable-name fix)
$ws.Cells.Item(2, 10).Value = [double]"20594.15739971409"
$ws.Cells.Item(2, 12).Value = [double]"3175.679150721244"
$ws.Cells.Item(2, 14).Value = [double]"20347.8347972425"
$ws.Cells.Item(2, 27).Value = [double]"0.7821998724435915"
$ws.Cells.Item(2, 29).Value = [double]"35.50917390327657"
$ws.Cells.Item(2, 30).Value = [double]"-4.547473508864641e-11"
$ws.Cells.Item(3, 9).Value = [double]"20594.15739971409"
$ws.Cells.Item(3, 10).Value = [double]"51605.72545261287"
$ws.Cells.Item(3, 11).Value = [double]"3175.679150721244"
$ws.Cells.Item(3, 12).Value = [double]"7957.753415052685"
$ws.Cells.Item(3, 13).Value = [double]"20347.8347972425"
$ws.Cells.Item(3, 14).Value = [double]"50988.47967997941"
$ws.Cells.Item(3, 26).Value = [double]"9.656032257748898"
$ws.Cells.Item(3, 27).Value = [double]"4.947756960897486"
$ws.Cells.Item(3, 28).Value = [double]"35.50917390327657"
$ws.Cells.Item(3, 29).Value = [double]"30.80089860642516"
$ws.Cells.Item(3, 30).Value = [double]"3.637978807091713e-11"
$ws.Cells.Item(4, 9).Value = [double]"51605.72545261287"
$ws.Cells.Item(4, 10).Value = [double]"85540.20770655214"
$ws.Cells.Item(4, 11).Value = [double]"7957.753415052685"
$ws.Cells.Item(4, 12).Value = [double]"13190.54957625184"
$ws.Cells.Item(4, 13).Value = [double]"50988.47967997941"
$ws.Cells.Item(4, 14).Value = [double]"84517.07837092172"
$ws.Cells.Item(4, 26).Value = [double]"12.0669967208405"
$ws.Cells.Item(4, 27).Value = [double]"8.534017753859334"
$ws.Cells.Item(4, 28).Value = [double]"30.80089860642516"
$ws.Cells.Item(4, 29).Value = [double]"27.26791963944399"
$ws.Cells.Item(4, 30).Value = [double]"-1.455191522836685e-10"
$ws.Cells.Item(5, 9).Value = [double]"85540.20770655214"
$ws.Cells.Item(5, 10).Value = [double]"118759.0138816608"
$ws.Cells.Item(5, 11).Value = [double]"13190.54957625184"
$ws.Cells.Item(5, 12).Value = [double]"18312.98639824133"
$ws.Cells.Item(5, 13).Value = [double]"84517.07837092172"
$ws.Cells.Item(5, 14).Value = [double]"117338.5610416385"
$ws.Cells.Item(5, 26).Value = [double]"14.24952565581575"
$ws.Cells.Item(5, 27).Value = [double]"11.6093100163225"
$ws.Cells.Item(5, 28).Value = [double]"27.26791963944399"
$ws.Cells.Item(5, 29).Value = [double]"24.62770399995075"
$ws.Cells.Item(5, 30).Value = [double]"0"
$ws.Cells.Item(6, 9).Value = [double]"118759.0138816608"
$ws.Cells.Item(6, 10).Value = [double]"149254.824326475"
$ws.Cells.Item(6, 11).Value = [double]"18312.98639824133"
$ws.Cells.Item(6, 12).Value = [double]"23015.52933477768"
$ws.Cells.Item(6, 13).Value = [double]"117338.5610416385"
$ws.Cells.Item(6, 14).Value = [double]"147469.6171900059"
$ws.Cells.Item(6, 26).Value = [double]"16.13714366130164"
$ws.Cells.Item(6, 27).Value = [double]"14.24457579670524"
$ws.Cells.Item(6, 28).Value = [double]"24.62770399995075"
$ws.Cells.Item(6, 29).Value = [double]"22.73513613535435"
$ws.Cells.Item(6, 30).Value = [double]"-6.257323548197746e-10"
$ws.Cells.Item(7, 9).Value = [double]"149254.824326475"
$ws.Cells.Item(7, 10).Value = [double]"175855.2450482356"
$ws.Cells.Item(7, 11).Value = [double]"23015.52933477768"
$ws.Cells.Item(7, 12).Value = [double]"27117.39181193256"
$ws.Cells.Item(7, 13).Value = [double]"147469.6171900059"
$ws.Cells.Item(7, 14).Value = [double]"173751.8755935979"
$ws.Cells.Item(7, 26).Value = [double]"17.72511809272836"
$ws.Cells.Item(7, 27).Value = [double]"16.49004237033597"
$ws.Cells.Item(7, 28).Value = [double]"22.73513613535435"
$ws.Cells.Item(7, 29).Value = [double]"21.50006041296196"
$ws.Cells.Item(7, 30).Value = [double]"-1.164153218269348e-10"
$ws.Cells.Item(8, 9).Value = [double]"175855.2450482356"
$ws.Cells.Item(8, 10).Value = [double]"199366.0711412441"
$ws.Cells.Item(8, 11).Value = [double]"27117.39181193256"
$ws.Cells.Item(8, 12).Value = [double]"30742.82978400695"
$ws.Cells.Item(8, 13).Value = [double]"173751.8755935979"
$ws.Cells.Item(8, 14).Value = [double]"196981.4933925702"
$ws.Cells.Item(8, 15).Value = [double]"0"
$ws.Cells.Item(8, 26).Value = [double]"19.00012660679674"
$ws.Cells.Item(8, 27).Value = [double]"18.20857042225747"
$ws.Cells.Item(8, 28).Value = [double]"21.50006041296196"
$ws.Cells.Item(8, 29).Value = [double]"20.70850422842268"
$ws.Cells.Item(8, 30).Value = [double]"-8.731149137020111e-11"
$ws.Cells.Item(9, 9).Value = [double]"199366.0711412441"
$ws.Cells.Item(9, 10).Value = [double]"220179.3636843024"
$ws.Cells.Item(9, 11).Value = [double]"30742.82978400695"
$ws.Cells.Item(9, 12).Value = [double]"33952.30021311853"
$ws.Cells.Item(9, 13).Value = [double]"196981.4933925702"
$ws.Cells.Item(9, 14).Value = [double]"217545.842301485"
$ws.Cells.Item(9, 15).Value = [double]"0"
$ws.Cells.Item(9, 26).Value = [double]"19.82924129800264"
$ws.Cells.Item(9, 27).Value = [double]"19.42191199952205"
$ws.Cells.Item(9, 28).Value = [double]"20.70850422842268"
$ws.Cells.Item(9, 29).Value = [double]"20.30117492994209"
$ws.Cells.Item(9, 30).Value = [double]"3.783497959375381e-10"
$ws.Cells.Item(10, 9).Value = [double]"220179.3636843024"
$ws.Cells.Item(10, 10).Value = [double]"236963.9727571273"
$ws.Cells.Item(10, 11).Value = [double]"33952.30021311853"
$ws.Cells.Item(10, 12).Value = [double]"36540.5358981734"
$ws.Cells.Item(10, 13).Value = [double]"217545.842301485"
$ws.Cells.Item(10, 14).Value = [double]"234129.6940183258"
$ws.Cells.Item(10, 15).Value = [double]"0"
$ws.Cells.Item(10, 26).Value = [double]"20.21169249292412"
$ws.Cells.Item(10, 27).Value = [double]"20.26541035623375"
$ws.Cells.Item(10, 28).Value = [double]"20.30117492994209"
$ws.Cells.Item(10, 29).Value = [double]"20.35489279325171"
$ws.Cells.Item(10, 30).Value = [double]"3.783497959375381e-10"
$ws.Cells.Item(11, 9).Value = [double]"236963.9727571273"
$ws.Cells.Item(11, 10).Value = [double]"250039.058231366"
$ws.Cells.Item(11, 11).Value = [double]"36540.5358981734"
$ws.Cells.Item(11, 12).Value = [double]"38556.75222247004"
$ws.Cells.Item(11, 13).Value = [double]"234129.6940183258"
$ws.Cells.Item(11, 14).Value = [double]"247048.3910072751"
$ws.Cells.Item(11, 26).Value = [double]"20.18007535100958"
$ws.Cells.Item(11, 27).Value = [double]"21.00187983491386"
$ws.Cells.Item(11, 28).Value = [double]"20.35489279325171"
$ws.Cells.Item(11, 29).Value = [double]"21.17669727715599"
$ws.Cells.Item(11, 30).Value = [double]"4.074536263942719e-10"
$ws.Cells.Item(12, 9).Value = [double]"250039.058231366"
$ws.Cells.Item(12, 10).Value = [double]"252038.8353915245"
$ws.Cells.Item(12, 11).Value = [double]"38556.75222247004"
$ws.Cells.Item(12, 12).Value = [double]"38865.12369455038"
$ws.Cells.Item(12, 13).Value = [double]"247048.3910072751"
$ws.Cells.Item(12, 14).Value = [double]"249024.2492323254"
$ws.Cells.Item(12, 26).Value = [double]"19.69091319659749"
$ws.Cells.Item(12, 27).Value = [double]"21.45449097001999"
$ws.Cells.Item(12, 28).Value = [double]"21.17669727715599"
$ws.Cells.Item(12, 29).Value = [double]"22.94027505057849"
$ws.Cells.Item(12, 30).Value = [double]"1.455191522836685e-10"
$ws.Cells.Item(13, 9).Value = [double]"252038.8353915245"
$ws.Cells.Item(13, 10).Value = [double]"242791.7302298158"
$ws.Cells.Item(13, 11).Value = [double]"38865.12369455038"
$ws.Cells.Item(13, 12).Value = [double]"37439.19310187788"
$ws.Cells.Item(13, 13).Value = [double]"249024.2492323254"
$ws.Cells.Item(13, 14).Value = [double]"239887.7468481206"
$ws.Cells.Item(13, 26).Value = [double]"19.03356116155288"
$ws.Cells.Item(13, 27).Value = [double]"21.46665524725301"
$ws.Cells.Item(13, 28).Value = [double]"22.94027505057849"
$ws.Cells.Item(13, 29).Value = [double]"25.37336913627862"
$ws.Cells.Item(13, 30).Value = [double]"-2.619344741106033e-10"
$ws.Cells.Item(14, 9).Value = [double]"242791.7302298158"
$ws.Cells.Item(14, 10).Value = [double]"223605.6789103798"
$ws.Cells.Item(14, 11).Value = [double]"37439.19310187788"
$ws.Cells.Item(14, 12).Value = [double]"34480.6480166273"
$ws.Cells.Item(14, 13).Value = [double]"239887.7468481206"
$ws.Cells.Item(14, 14).Value = [double]"220931.1760556335"
$ws.Cells.Item(14, 26).Value = [double]"18.11164607798419"
$ws.Cells.Item(14, 27).Value = [double]"21.4138440438239"
$ws.Cells.Item(14, 28).Value = [double]"25.37336913627862"
$ws.Cells.Item(14, 29).Value = [double]"28.67556710211832"
$ws.Cells.Item(14, 30).Value = [double]"3.783497959375381e-10"
$ws.Cells.Item(15, 9).Value = [double]"223605.6789103798"
$ws.Cells.Item(15, 10).Value = [double]"195563.0847811293"
$ws.Cells.Item(15, 11).Value = [double]"34480.6480166273"
$ws.Cells.Item(15, 12).Value = [double]"30156.39819276051"
$ws.Cells.Item(15, 13).Value = [double]"220931.1760556335"
$ws.Cells.Item(15, 14).Value = [double]"193223.9937925693"
$ws.Cells.Item(15, 26).Value = [double]"16.89513507571947"
$ws.Cells.Item(15, 27).Value = [double]"21.11422971377652"
$ws.Cells.Item(15, 28).Value = [double]"28.67556710211832"
$ws.Cells.Item(15, 29).Value = [double]"32.89466174017537"
$ws.Cells.Item(15, 30).Value = [double]"8.731149137020111e-11"
$ws.Cells.Item(16, 9).Value = [double]"195563.0847811293"
$ws.Cells.Item(16, 10).Value = [double]"159986.0091855307"
$ws.Cells.Item(16, 11).Value = [double]"30156.39819276051"
$ws.Cells.Item(16, 12).Value = [double]"24670.30934631202"
$ws.Cells.Item(16, 13).Value = [double]"193223.9937925693"
$ws.Cells.Item(16, 14).Value = [double]"158072.448490779"
$ws.Cells.Item(16, 26).Value = [double]"15.27785613197528"
$ws.Cells.Item(16, 27).Value = [double]"20.53791586463613"
$ws.Cells.Item(16, 28).Value = [double]"32.89466174017537"
$ws.Cells.Item(16, 29).Value = [double]"38.15472147283623"
$ws.Cells.Item(16, 30).Value = [double]"5.529727786779404e-10"
$ws.Cells.Item(17, 9).Value = [double]"159986.0091855307"
$ws.Cells.Item(17, 10).Value = [double]"119026.9425532871"
$ws.Cells.Item(17, 11).Value = [double]"24670.30934631202"
$ws.Cells.Item(17, 12).Value = [double]"18354.30178103897"
$ws.Cells.Item(17, 13).Value = [double]"158072.448490779"
$ws.Cells.Item(17, 14).Value = [double]"117603.2850719488"
$ws.Cells.Item(17, 26).Value = [double]"13.14129890466775"
$ws.Cells.Item(17, 27).Value = [double]"19.62809090534333"
$ws.Cells.Item(17, 28).Value = [double]"38.15472147283623"
$ws.Cells.Item(17, 29).Value = [double]"44.64151347351181"
$ws.Cells.Item(17, 30).Value = [double]"-2.328306436538696e-10"
$ws.Cells.Item(18, 9).Value = [double]"119026.9425532871"
$ws.Cells.Item(18, 10).Value = [double]"75619.5729149049"
$ws.Cells.Item(18, 11).Value = [double]"18354.30178103897"
$ws.Cells.Item(18, 12).Value = [double]"11660.75874974338"
$ws.Cells.Item(18, 13).Value = [double]"117603.2850719488"
$ws.Cells.Item(18, 14).Value = [double]"74715.10231012804"
$ws.Cells.Item(18, 26).Value = [double]"10.26784755027579"
$ws.Cells.Item(18, 27).Value = [double]"18.31918012529943"
$ws.Cells.Item(18, 28).Value = [double]"44.64151347351181"
$ws.Cells.Item(18, 29).Value = [double]"52.69284604853544"
$ws.Cells.Item(18, 30).Value = [double]"1.673470251262188e-10"
$ws.Cells.Item(19, 9).Value = [double]"75619.5729149049"
$ws.Cells.Item(19, 10).Value = [double]"42222.30508670387"
$ws.Cells.Item(19, 11).Value = [double]"11660.75874974338"
$ws.Cells.Item(19, 12).Value = [double]"6510.802620217835"
$ws.Cells.Item(19, 13).Value = [double]"74715.10231012804"
$ws.Cells.Item(19, 14).Value = [double]"41717.29253049944"
$ws.Cells.Item(19, 26).Value = [double]"7.772449164699324"
$ws.Cells.Item(19, 27).Value = [double]"15.73623342099551"
$ws.Cells.Item(19, 28).Value = [double]"52.69284604853544"
$ws.Cells.Item(19, 29).Value = [double]"60.65663030483163"
$ws.Cells.Item(19, 30).Value = [double]"7.275957614183426e-12"
$ws.Cells.Item(20, 9).Value = [double]"42222.30508670387"
$ws.Cells.Item(20, 10).Value = [double]"13221.40609254259"
$ws.Cells.Item(20, 11).Value = [double]"6510.802620217835"
$ws.Cells.Item(20, 12).Value = [double]"2038.779390502728"
$ws.Cells.Item(20, 13).Value = [double]"41717.29253049944"
$ws.Cells.Item(20, 14).Value = [double]"13063.26749556877"
$ws.Cells.Item(20, 26).Value = [double]"3.953103522370313"
$ws.Cells.Item(20, 27).Value = [double]"15.35253924342194"
$ws.Cells.Item(20, 28).Value = [double]"60.65663030483163"
$ws.Cells.Item(20, 29).Value = [double]"72.05606602588325"
$ws.Cells.Item(20, 30).Value = [double]"-2.546585164964199e-11"
$ws.Cells.Item(21, 9).Value = [double]"13221.40609254259"
$ws.Cells.Item(21, 11).Value = [double]"2038.779390502728"
$ws.Cells.Item(21, 13).Value = [double]"13063.26749556877"
$ws.Cells.Item(21, 26).Value = [double]"3.953103522370313"
$ws.Cells.Item(21, 27).Value = [double]"15.35253924342194"
$ws.Cells.Item(21, 28).Value = [double]"72.05606602588325"

Write-Output "applied 2 header + 213 data changes"
